# prevention and rehab 2.0.0 added
$wb = $excel.ActiveWorkbook

# --- Remove the "Include from Tempcodes" sheet -----------------------------
# Its single concept row gets folded into "Include from CareSocialCodes"
# below, so the whole worksheet goes away.
$tempSheet = $wb.Worksheets.Item("Include from Tempcodes")
$tempSheet.Delete() | Out-Null

# --- Metadata sheet: bump version / date / contact --------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-14T10:48:54+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- Include from CareSocialCodes: merge in the Tempcodes concept ----------
$care = $wb.Worksheets.Item("Include from CareSocialCodes")

# Row 10 used to be a blank placeholder row (empty concept + empty
# description); it now carries the UUID that used to live on the
# Tempcodes sheet, with an empty description cell.
$care.Range("A10").Value = "ef491570-7884-4acd-bcf7-43d6b2c64762"
$care.Range("B10").Value = ""

# Row 11 used to hold the "System URI" / CareSocialCodes URL pair; that pair
# moves down to the new row 12, so row 11 becomes blank (empty text, same
# as the other placeholder rows above it).
$care.Range("A11").Value = "'"
$care.Range("B11").Value = "'"
$care.Range("A9:B9").Copy()
$care.Range("A11:B11").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# New row 12: the "System URI" row, re-added at the bottom with the same
# formatting as the other data rows.
$care.Range("A9:B9").Copy()
$care.Range("A12:B12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$care.Range("A12").Value = "System URI"
$care.Range("B12").Value = "http://fhir.kl.dk/term/CodeSystem/CareSocialCodes"
